$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.699.99'
$ws.Range("E2").Value = '  -0.26%  '
$ws.Range("D3").Value = '1.629.10'
$ws.Range("E3").Value = '  -0.30%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = '  -0.77%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("E8").Value = '  -1.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0631'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.34%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.44'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.81%  '
$ws.Range("E11").Value = '  +1.70%  '
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("D13").Value = '1.854.47'
$ws.Range("E13").Value = '  -0.26%  '
$ws.Range("D14").Value = '1.644.46'
$ws.Range("E14").Value = '  +0.57%  '
$ws.Range("E15").Value = '  -0.05%  '
$ws.Range("E16").Value = '  -2.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.80'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.42%  '
$ws.Range("D18").Value = '25.713.57'
$ws.Range("E18").Value = '  -0.27%  '
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '191.48'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.30%  '
$ws.Range("E22").Value = '  -0.55%  '
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("E25").Value = '  +2.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.39'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.95%  '
$ws.Range("E27").Value = '  +1.93%  '
$ws.Range("E28").Value = '  +0.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.44'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.81%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.23'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0489'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.30%  '
$ws.Range("E32").Value = '  -0.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.22'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.03%  '
$ws.Range("E34").Value = '  -0.74%  '
$ws.Range("E35").Value = '  +0.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.902'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.67%  '
$ws.Range("D37").Value = '1.136.79'
$ws.Range("E37").Value = '  +2.71%  '
$ws.Range("E38").Value = '  -2.32%  '
$ws.Range("E40").Value = '  -0.91%  '
$ws.Range("B41").Value = 'mCoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.54'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.26%  '
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.37'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.18%  '
$ws.Range("E44").Value = '  -0.58%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.801'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.06%  '
$ws.Range("D46").Value = '1.763.50'
$ws.Range("E46").Value = '  -0.09%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.0₆0109'
$ws.Range("E47").Value = '  +1.10%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.15'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("E49").Value = '  +0.15%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0507'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.56%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.43'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.15%  '
